# Make sure "smart quotes" autocorrect doesn't mangle straight apostrophes
# while we edit (belt-and-braces; the XML inserts below are immune anyway).
$word.Options.AutoFormatReplaceQuotes = $false
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Remove the "Meta description: ..." paragraph that follows the title ---
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# --- 2. Insert a new bold "Play Fruits & 777's Free..." paragraph right  ---
#        before the final ("Create a feature image...") paragraph.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$titleParaXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruits &amp; 777's Free: A Classic Fruit Slot Game</w:t></w:r></w:p><w:p $wNs/>"
$null = $insertionPoint.InsertXML($titleParaXml)

# InsertXML above adds both the real paragraph and a trailing blank
# placeholder paragraph (needed to force the paragraph break); drop the
# placeholder now that it has served its purpose.
$placeholder = $d.Paragraphs($count + 1)
$placeholder.Range.Delete()

# --- 3. Replace the old "Create a feature image..." prompt text with the ---
#        meta-description copy (keeping its existing italic formatting).
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$descParaXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Check out our review of Fruits &amp; 777's from Spearhead Studios, a classic fruit slot game with simple gameplay, reasonable RTP range, and a wide betting range to please all players. Play for free now!</w:t></w:r></w:p>"
$null = $finalPara.Range.InsertXML($descParaXml)
